$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 ("Number of registered, unit") updated values E4:K4
$ws.Range("E4").Value = 4182
$ws.Range("F4").Value = 3652
$ws.Range("G4").Value = 3446
$ws.Range("H4").Value = 3424
$ws.Range("I4").Value = 3358
$ws.Range("J4").Value = 3434
$ws.Range("K4").Value = 3471

# Row 5 ("Number of beneficiaries, unit") updated values E5:K5
$ws.Range("E5").Value = 1937
$ws.Range("F5").Value = 1956
$ws.Range("G5").Value = 1820
$ws.Range("H5").Value = 1964
$ws.Range("I5").Value = 1701
$ws.Range("J5").Value = 2030
$ws.Range("K5").Value = 2110

# Row 5 switches to the same (right-aligned) number format used by F4:K4,
# replacing its previous left-aligned/bordered style.
$ws.Range("F4").Copy() | Out-Null
$ws.Range("E5:K5").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false | Out-Null

# Move/save the active selection to A3, matching the updated view state
$ws.Range("A3").Select() | Out-Null
